# Sprint 3 burndown update
# - Fill in actual hours for the "Class Demo" task (row 8) on Sprint3
# - Fill in actual burndown hours achieved (G19, G26, G27) on Sprint3 so the
#   burndown totals/chart (G35:J36) recompute
# - Leave the workbook with ProductBacklog as the active sheet/selection,
#   matching the author's final view state

$wb = $excel.ActiveWorkbook

$sprint3 = $wb.Worksheets.Item("Sprint3")

# "Class Demo" task hours (C8:E8) -> totals/formulas recompute automatically
$sprint3.Range("C8").Value = 6
$sprint3.Range("D8").Value = 6
$sprint3.Range("E8").Value = 6

# Actual hours completed, feeding the burndown chart (G35:J36 recompute)
$sprint3.Range("G19").Value = 2
$sprint3.Range("G26").Value = 3
$sprint3.Range("G27").Value = 3

# Leave Sprint2 view at its working position (was previously the selected tab)
$sprint2 = $wb.Worksheets.Item("Sprint2")
$sprint2.Activate() | Out-Null
$sprint2.Range("J21").Select() | Out-Null

# Leave Sprint3 view selected near the new data
$sprint3.Activate() | Out-Null
$sprint3.Range("M24").Select() | Out-Null

# Final state: ProductBacklog is the active/selected sheet, cursor on A6
$backlog = $wb.Worksheets.Item("ProductBacklog")
$backlog.Activate() | Out-Null
$backlog.Range("A6").Select() | Out-Null

Write-Output "Sprint 3 burndown updated"
